# The sheet originally tracked 6 "Endorsement" columns (B:G) in between the
# "Instrument Type*" column (A) and the "Company Master Id*" column (H).
# Those endorsement-tracking columns are no longer needed, so remove them -
# this shifts "Company Master Id*" left into column B and shrinks the used
# range from A1:H10 down to A1:B10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colsToRemove = 6
$ws.Range("B1:G1").EntireColumn.Delete()

# Deleting whole columns shrinks the sheet's trailing "default width" column
# stretch (normally anchored at column 16384/XFD) by the number of columns
# removed. Re-insert the same number of blank columns at the very end of the
# sheet (well past any real data) so the column-width metadata once again
# spans out to the worksheet's true last column, matching stock Excel
# behaviour where the address space stays fixed at 16384 columns.
$totalCols = $ws.Columns.Count
for ($i = 0; $i -lt $colsToRemove; $i++) {
    $ws.Columns.Item($totalCols - $colsToRemove + $i).EntireColumn.Insert()
}
